# Updated ExperimentDescription with results, and added some more notes to results sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / add the note texts on the results sheet ---
# Editing A41 first (in place) reuses the existing shared-string slot for the old
# note text, turning it into "Note 2: ...". Then two brand-new notes are entered
# ("Note 3: ..." which re-uses the original note wording, and "Note 1: ...").
# Doing this in this order keeps the shared string table ordering consistent with
# the target workbook: Note 2 keeps the old index, Note 3 is appended next, and
# Note 1 is appended last.
$ws.Range("A41").Value = "Note 2: values of 0.000 were entered as 0.0005, so that they could still be graphed on the log10 scale. Tests above 60 real-time seconds were aborted automatically by the failsafe."
$ws.Range("A11").Value = "Note 3: For the Vector Write Test, the 1-minute cap was disabled, as it seemed likely it would still finish within ~100 seconds, and its value was still reasonable to be recorded. The Read and Write Test times for Linked Lists, for comparison, were getting absurd."
$ws.Range("A9").Value = "Note 1: For the linked list read and write tests, times were heavily approximated. See ExperimentDescription.txt for details."

# Move the "Note 2" text up from row 41 into row 10, right under Note 1, and
# remove the now empty row 41 so the note block sits directly under the data
# (rows 9-11) instead of way down at row 41.
$ws.Range("A10").Value = $ws.Range("A41").Value()
$ws.Rows.Item(41).Delete()

# --- Reposition the three charts slightly (they were nudged down/left a bit) ---
$charts = $ws.ChartObjects()

$c1 = $charts.Item(1)
$c1.Top = 169.2
$c1.Left = 4.8

$c2 = $charts.Item(2)
$c2.Top = 166.2
$c2.Left = 474.1

$c3 = $charts.Item(3)
$c3.Top = 404.4
$c3.Left = 222.1125

# --- Update the active selection shown when the workbook is reopened ---
$ws.Range("S9").Select()
